$d = $word.ActiveDocument

$replacements = @(
    @{old = "296×4=1184"; new = "627×2=1254"},
    @{old = "235×7=1645"; new = "980×9=8820"},
    @{old = "940×9=8460"; new = "698×2=1396"},
    @{old = "205×3=615";  new = "594×8=4752"},
    @{old = "588×4=2352"; new = "810×8=6480"},
    @{old = "292×4=1168"; new = "155×6=930"},
    @{old = "518×9=4662"; new = "712×6=4272"},
    @{old = "318×4=1272"; new = "469×7=3283"},
    @{old = "131×5=655";  new = "856×8=6848"},
    @{old = "746×8=5968"; new = "183×5=915"},
    @{old = "464×2=928";  new = "595×7=4165"},
    @{old = "585×8=4680"; new = "392×8=3136"},
    @{old = "193×5=965";  new = "481×3=1443"},
    @{old = "795×6=4770"; new = "168×4=672"},
    @{old = "873×4=3492"; new = "612×8=4896"},
    @{old = "386×5=1930"; new = "960×9=8640"},
    @{old = "438×9=3942"; new = "289×9=2601"},
    @{old = "828×4=3312"; new = "488×6=2928"},
    @{old = "489×7=3423"; new = "590×4=2360"},
    @{old = "759×4=3036"; new = "367×9=3303"},
    @{old = "384×9=3456"; new = "714×9=6426"},
    @{old = "153×8=1224"; new = "453×7=3171"},
    @{old = "507×7=3549"; new = "683×5=3415"},
    @{old = "505×7=3535"; new = "627×4=2508"},
    @{old = "464×5=2320"; new = "555×3=1665"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
